$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "67.369.77"
$ws.Range("D3").Value = "2.553.19"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.61"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.07"
$ws.Range("E6").Value = "  +4.97%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "2.553.38"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.17"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.14"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "3.004.80"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "67.176.41"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "2.542.10"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.45"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.81"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +6.21%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.23"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.74"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("E32").Value = "  +5.27%  "
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.134"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.38"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.47"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.22"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("E44").Value = "  +6.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.84"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.44"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.567"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0284"
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("E51").Value = "  +1.56%  "
